# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals.
# Only column G ("K") values for rows 2-12 change.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 1
    5  = 0
    6  = 2
    7  = 2
    8  = 0
    9  = 1
    10 = 1
    11 = 3
    12 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
